{"js": "const replacements = [\n    [\"2024-09-05 Thursday\", \"2024-09-06 Friday\"],\n    [\"36\u00f75=\", \"25\u00f76=\"],\n    [\"13\u00f77=\", \"95\u00f78=\"],\n    [\"73\u00f73=\", \"59\u00f77=\"],\n    [\"76\u00f72=\", \"23\u00f79=\"],\n    [\"64\u00f79=\", \"53\u00f79=\"],\n    [\"52\u00f73=\", \"85\u00f78=\"],\n    [\"76\u00f73=\", \"56\u00f78=\"],\n    [\"57\u00f77=\", \"71\u00f77=\"],\n    [\"67\u00f75=\", \"10\u00f74=\"],\n    [\"26\u00f77=\", \"52\u00f76=\"],\n    [\"52\u00f72=\", \"72\u00f79=\"],\n    [\"32\u00f78=\", \"98\u00f79=\"],\n    [\"80\u00f74=\", \"21\u00f73=\"],\n    [\"19\u00f75=\", \"24\u00f78=\"],\n    [\"90\u00f77=\", \"31\u00f78=\"],\n    [\"32\u00f72=\", \"58\u00f76=\"],\n    [\"12\u00f73=\", \"76\u00f75=\"],\n    [\"99\u00f79=\", \"94\u00f75=\"],\n    [\"91\u00f72=\", \"92\u00f79=\"],\n    [\"97\u00f79=\", \"30\u00f78=\"],\n    [\"57\u00f79=\", \"35\u00f77=\"],\n    [\"70\u00f76=\", \"21\u00f78=\"],\n    [\"22\u00f75=\", \"52\u00f79=\"],\n    [\"86\u00f72=\", \"36\u00f79=\"],\n    [\"45\u00f76=\", \"34\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2024-09-05 Thursday\"; New = \"2024-09-06 Friday\" },\n    @{ Old = \"36\u00f75=\"; New = \"25\u00f76=\" },\n    @{ Old = \"13\u00f77=\"; New = \"95\u00f78=\" },\n    @{ Old = \"73\u00f73=\"; New = \"59\u00f77=\" },\n    @{ Old = \"76\u00f72=\"; New = \"23\u00f79=\" },\n    @{ Old = \"64\u00f79=\"; New = \"53\u00f79=\" },\n    @{ Old = \"52\u00f73=\"; New = \"85\u00f78=\" },\n    @{ Old = \"76\u00f73=\"; New = \"56\u00f78=\" },\n    @{ Old = \"57\u00f77=\"; New = \"71\u00f77=\" },\n    @{ Old = \"67\u00f75=\"; New = \"10\u00f74=\" },\n    @{ Old = \"26\u00f77=\"; New = \"52\u00f76=\" },\n    @{ Old = \"52\u00f72=\"; New = \"72\u00f79=\" },\n    @{ Old = \"32\u00f78=\"; New = \"98\u00f79=\" },\n    @{ Old = \"80\u00f74=\"; New = \"21\u00f73=\" },\n    @{ Old = \"19\u00f75=\"; New = \"24\u00f78=\" },\n    @{ Old = \"90\u00f77=\"; New = \"31\u00f78=\" },\n    @{ Old = \"32\u00f72=\"; New = \"58\u00f76=\" },\n    @{ Old = \"12\u00f73=\"; New = \"76\u00f75=\" },\n    @{ Old = \"99\u00f79=\"; New = \"94\u00f75=\" },\n    @{ Old = \"91\u00f72=\"; New = \"92\u00f79=\" },\n    @{ Old = \"97\u00f79=\"; New = \"30\u00f78=\" },\n    @{ Old = \"57\u00f79=\"; New = \"35\u00f77=\" },\n    @{ Old = \"70\u00f76=\"; New = \"21\u00f78=\" },\n    @{ Old = \"22\u00f75=\"; New = \"52\u00f79=\" },\n    @{ Old = \"86\u00f72=\"; New = \"36\u00f79=\" },\n    @{ Old = \"45\u00f76=\"; New = \"34\u00f76=\" }\n)\n\nforeach ($pair in $replacements) {\n    $r = $d.Content\n    $r.Find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
